$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3020.0557
$ws.Range("J70").Value = 3157.8333
$ws.Range("L70").Value = 9473.499899999999
$ws.Range("N70").Value = -10013.4999
$ws.Range("H73").Value = 3020.0557
$ws.Range("J73").Value = 3157.8333
$ws.Range("L73").Value = 9473.499899999999
$ws.Range("N73").Value = -11345.4999
$ws.Range("H76").Value = 6327.0527
$ws.Range("J76").Value = 6791.4614
$ws.Range("L76").Value = 6791.4614
$ws.Range("N76").Value = -7421.4614
$ws.Range("H79").Value = 6327.0527
$ws.Range("J79").Value = 6791.4614
$ws.Range("L79").Value = 6791.4614
$ws.Range("N79").Value = -8975.4614
$ws.Range("H116").Value = 2998.35
$ws.Range("I116").Value = 3170
$ws.Range("K116").Value = 3170
$ws.Range("M116").Value = 272
$ws.Range("H132").Value = 4250.1665
$ws.Range("I132").Value = 4609.3335
$ws.Range("K132").Value = 13828.0005
$ws.Range("M132").Value = -11298.0005
$ws.Range("H136").Value = 109999.5
$ws.Range("J136").Value = 109999.5
$ws.Range("L136").Value = 109999.5
$ws.Range("N136").Value = -120199.5
$ws.Range("H137").Value = 3120.8125
$ws.Range("I137").Value = 3256.6667
$ws.Range("K137").Value = 9770.000100000001
$ws.Range("M137").Value = -7220.000100000001
$ws.Range("H138").Value = 5295308
$ws.Range("I138").Value = 2514.5715
$ws.Range("J138").Value = 6807534.5
$ws.Range("K138").Value = 7543.7145
$ws.Range("L138").Value = 20422603.5
$ws.Range("M138").Value = -2403.7145
$ws.Range("N138").Value = -20432883.5
$ws.Range("H140").Value = 99999.5
$ws.Range("J140").Value = 99999.5
$ws.Range("L140").Value = 99999.5
$ws.Range("N140").Value = -110359.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15395800
$ws.Range("I32").Value = 24395132
$ws.Range("K32").Value = 24395132
$ws.Range("M32").Value = -24394845
$ws.Range("H61").Value = 29415372
$ws.Range("I61").Value = 47621550
$ws.Range("J61").Value = 5388.4614
$ws.Range("K61").Value = 47621550
$ws.Range("L61").Value = 5388.4614
$ws.Range("M61").Value = -47621338
$ws.Range("N61").Value = -5812.4614
$ws.Range("H102").Value = 106753
$ws.Range("I102").Value = 112350.39
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 112350.39
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -110728.39
$ws.Range("N102").Value = -9244
$ws.Range("H122").Value = 3644.2222
$ws.Range("I122").Value = 1872.7142
$ws.Range("K122").Value = 5618.142599999999
$ws.Range("M122").Value = -3168.142599999999
$ws.Range("H132").Value = 28579750
$ws.Range("I132").Value = 9723.071
$ws.Range("J132").Value = 142859860
$ws.Range("K132").Value = 29169.213
$ws.Range("L132").Value = 428579580
$ws.Range("M132").Value = -26639.213
$ws.Range("N132").Value = -428584640
$ws.Range("H136").Value = 29415372
$ws.Range("I136").Value = 47621550
$ws.Range("J136").Value = 5388.4614
$ws.Range("K136").Value = 142864650
$ws.Range("L136").Value = 16165.3842
$ws.Range("M136").Value = -142862100
$ws.Range("N136").Value = -21265.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 21007.334
$ws.Range("I44").Value = 21248.5
$ws.Range("K44").Value = 21248.5
$ws.Range("M44").Value = -20751.5
$ws.Range("H99").Value = 3333.4167
$ws.Range("I99").Value = 2475.1
$ws.Range("J99").Value = 3946.5
$ws.Range("K99").Value = 2475.1
$ws.Range("L99").Value = 3946.5
$ws.Range("M99").Value = -977.0999999999999
$ws.Range("N99").Value = -6942.5
$ws.Range("H134").Value = 5126.1035
$ws.Range("I134").Value = 4826.1665
$ws.Range("K134").Value = 14478.4995
$ws.Range("M134").Value = -11943.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1687
$ws.Range("I122").Value = 1714.1818
$ws.Range("J122").Value = 1612.25
$ws.Range("K122").Value = 5142.5454
$ws.Range("L122").Value = 4836.75
$ws.Range("M122").Value = -2692.5454
$ws.Range("N122").Value = -9736.75
$ws.Range("H131").Value = 71239
$ws.Range("J131").Value = 71239
$ws.Range("L131").Value = 71239
$ws.Range("N131").Value = -81319
$ws.Range("H134").Value = 2433.3076
$ws.Range("I134").Value = 1990.375
$ws.Range("J134").Value = 3142
$ws.Range("K134").Value = 5971.125
$ws.Range("L134").Value = 9426
$ws.Range("M134").Value = -3436.125
$ws.Range("N134").Value = -14496
$ws.Range("H141").Value = 317236.12
$ws.Range("J141").Value = 345412.72
$ws.Range("L141").Value = 345412.72
$ws.Range("N141").Value = -355772.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1402611.4
$ws.Range("I2").Value = 928.8889
$ws.Range("J2").Value = 2664125.8
$ws.Range("K2").Value = 5573.3334
$ws.Range("L2").Value = 15984754.8
$ws.Range("M2").Value = -5460.3334
$ws.Range("N2").Value = -15984980.8
$ws.Range("H113").Value = 2165.742
$ws.Range("I113").Value = 1674.6666
$ws.Range("J113").Value = 2366.6365
$ws.Range("K113").Value = 5023.9998
$ws.Range("L113").Value = 7099.9095
$ws.Range("M113").Value = -2853.9998
$ws.Range("N113").Value = -11439.9095
$ws.Range("H128").Value = 116092.5
$ws.Range("I128").Value = 116092.5
$ws.Range("K128").Value = 348277.5
$ws.Range("M128").Value = -343297.5
$ws.Range("H136").Value = 7391.2856
$ws.Range("J136").Value = 9916.666999999999
$ws.Range("L136").Value = 29750.001
$ws.Range("N136").Value = -39950.001
$ws.Range("H137").Value = 2878.077
$ws.Range("I137").Value = 1857
$ws.Range("J137").Value = 3753.2856
$ws.Range("K137").Value = 5571
$ws.Range("L137").Value = 11259.8568
$ws.Range("M137").Value = -471
$ws.Range("N137").Value = -21459.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2025.919
$ws.Range("J102").Value = 4061.1
$ws.Range("L102").Value = 4061.1
$ws.Range("N102").Value = -7305.1
$ws.Range("H132").Value = 2114.0588
$ws.Range("I132").Value = 1730.3077
$ws.Range("J132").Value = 3361.25
$ws.Range("K132").Value = 5190.9231
$ws.Range("L132").Value = 10083.75
$ws.Range("M132").Value = -2660.9231
$ws.Range("N132").Value = -15143.75
$ws.Range("H136").Value = 58078.25
$ws.Range("J136").Value = 58078.25
$ws.Range("L136").Value = 174234.75
$ws.Range("N136").Value = -179334.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2648.0908
$ws.Range("I22").Value = 1912.0834
$ws.Range("J22").Value = 3068.6667
$ws.Range("K22").Value = 1912.0834
$ws.Range("L22").Value = 3068.6667
$ws.Range("M22").Value = -1617.0834
$ws.Range("N22").Value = -3658.6667
$ws.Range("H27").Value = 2648.0908
$ws.Range("I27").Value = 1912.0834
$ws.Range("J27").Value = 3068.6667
$ws.Range("K27").Value = 1912.0834
$ws.Range("L27").Value = 3068.6667
$ws.Range("M27").Value = -1805.0834
$ws.Range("N27").Value = -3282.6667
$ws.Range("H122").Value = 4165.086
$ws.Range("I122").Value = 3771.12
$ws.Range("J122").Value = 5150
$ws.Range("K122").Value = 11313.36
$ws.Range("L122").Value = 15450
$ws.Range("M122").Value = -8863.360000000001
$ws.Range("N122").Value = -20350
$ws.Range("H136").Value = 4913.171
$ws.Range("J136").Value = 8906.571
$ws.Range("L136").Value = 26719.713
$ws.Range("N136").Value = -31819.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 20000006
$ws.Range("J10").Value = 20000006
$ws.Range("L10").Value = 20000006
$ws.Range("N10").Value = -20000344
$ws.Range("H58").Value = 31900
$ws.Range("I58").Value = 31900
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 31900
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -31592
$ws.Range("N58").ClearContents()
$ws.Range("H114").Value = 78000
$ws.Range("J114").Value = 78000
$ws.Range("L114").Value = 78000
$ws.Range("N114").Value = -86678
$ws.Range("H132").Value = 4525.6787
$ws.Range("J132").Value = 4898.3335
$ws.Range("L132").Value = 14695.0005
$ws.Range("N132").Value = -19755.0005
